$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 285
$ws1.Range("F4").Value = 1091
$ws1.Range("F5").Value = 571

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 285
$ws4.Range("F4").Value = 1091
$ws4.Range("F6").Value = 571
